$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: shift existing value ("Cache size is 1000 requests...") from B9 to C9,
# and set B9 to "yes"
$ws.Range("C9").Value = $ws.Range("B9").Text
$ws.Range("B9").Value = "yes"

# Row 16: add new cells B16 and C16
$ws.Range("B16").Value = "Strong exception "
$ws.Range("C16").Value = "write in danger log"

# Update the view: scroll/zoom and selection
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.Zoom = 92
$ws.Range("C12").Select()
